# Update COVID country statistics workbook + re-sort by "Casos totales" desc.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Refresh the "last updated" timestamp.
$ws.Range("A1").Value = "Datos actualizados a 21 de Abril de 2020 a las 00:52"

# 2) Apply the new per-country figures (Casos totales, Nuevos casos, Casos
#    activos, Recuperados, Casos criticos, Muertes hoy, Muertes).
#    Rows below refer to the CURRENT (pre-sort) layout of the sheet.
$updates = @{
    4   = @(789383, 24747, 71832, 675248, 13634, 1728, 42303)  # Estados Unidos
    8   = @(146777, 1035, 91500, 50475, 2889, 160, 4802)       # Alemania
    38  = @(7156, 78, 32, 6943, 58, 16, 181)                   # Noruega
    39  = @(6900, 154, 1559, 5147, 84, 8, 194)                 # Chequia
    42  = @(6625, 13, 4258, 2296, 49, 0, 71)                   # Australia
    84  = @(1163, 146, 305, 816, 33, 0, 42)                    # Camerun
    96  = @(665, 38, 188, 455, 2, 1, 22)                       # Nigeria
    124 = @(246, 22, 56, 176, 0, 0, 14)                        # Mali
    156 = @(66, 5, 7, 58, 3, 0, 1)                             # Islas Caimanes
    203 = @(9, 1, 2, 7, 0, 0, 0)                                # Santa Sede
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    for ($col = 2; $col -le 8; $col++) {
        $ws.Cells.Item($row, $col).Value = $vals[$col - 2]
    }
}

# 3) Re-sort the whole country table (rows 4-216) by "Casos totales"
#    descending, exactly like the live leaderboard does on every refresh.
$dataRange = $ws.Range("A4:H216")
$sortKey = $ws.Range("B4:B216")
$dataRange.Sort($sortKey, 2)
